# Update attributions.xlsx - "Met à jour attributions.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates (column A = task description, column B = assignee) ---

# Row 5: rename "question" -> "sondage"
$ws.Range("A5").Value = "Le formateur lance un sondage (oui/non, à choix multiples, à saisie libre) [2]"

# Row 7: rename "question" -> "sondage"
$ws.Range("A7").Value = "Le participant affiche les réponses à un sondage [1]"

# Row 8: rename "questions du canal" -> "sondages du canal"
$ws.Range("A8").Value = "Le participant liste les sondages du canal [1]"

# Row 13: assign Ilyesse
$ws.Range("B13").Value = "Ilyesse"

# Row 15: "toutes sessions confondues" -> "tous canaux confondus"
$ws.Range("A15").Value = "Le formateur liste ses questionnaires (tous canaux confondus). [1]"

# Row 16: new shorter task text, and un-assign (was Dieynaba)
$ws.Range("A16").Value = "Le formateur crée un questionnaire. [2]"
$ws.Range("B16").ClearContents()

# Row 18: updated task text, now assigned to Dieynaba
$ws.Range("A18").Value = "L'étudiant répond à un questionnaire. [1]"
$ws.Range("B18").Value = "Dieynaba"

# Row 27: replaced task (sessions -> canaux)
$ws.Range("A27").Value = "Le membre peut lister les canaux. [1] "

# Row 30: replaced task text
$ws.Range("A30").Value = "Le gestionnaire ou le formateur liste les membres d'un canal [1]"

# Row 31: assign Ilyesse
$ws.Range("B31").Value = "Ilyesse"

# --- Row heights: rows 5, 11, 27, 30 lose their explicit 30pt height (back to default) ---
$ws.Cells.Item(5, 1).EntireRow.AutoFit()
$ws.Cells.Item(11, 1).EntireRow.AutoFit()
$ws.Cells.Item(27, 1).EntireRow.AutoFit()
$ws.Cells.Item(30, 1).EntireRow.AutoFit()

# --- Column A width: 77.85546875 -> 78.7109375 (closest reachable width) ---
$ws.Columns.Item(1).ColumnWidth = 77.8

# --- Selection moves from A19 to A10 ---
$ws.Range("A10").Select() | Out-Null
